# Weekly fruit/vegetable data update: insert a new record as row 367,
# pushing all existing rows from 367..389 down to 368..390.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 367 (this shifts rows
# 367-389 down to 368-390, and the sheet dimension grows to A1:R390).
$ws.Rows.Item(367).Insert()

# Copy the date style (s="2") used by the rest of column D onto the new
# row's D cell, so the new date renders the same way as its neighbours.
$ws.Cells.Item(368, 4).Copy()
$ws.Cells.Item(367, 4).PasteSpecial(-4122)

# Populate the brand-new row 367 with the latest weekly observation.
$ws.Cells.Item(367, 1).Value = 6
$ws.Cells.Item(367, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(367, 3).Value = "Metropolitana"
$ws.Cells.Item(367, 4).Value = 45267
$ws.Cells.Item(367, 5).Value = 13
$ws.Cells.Item(367, 6).Value = 100112001
$ws.Cells.Item(367, 7).Value = "Berenjena"
$ws.Cells.Item(367, 8).Value = "Sin especificar"
$ws.Cells.Item(367, 9).Value = "Primera"
$ws.Cells.Item(367, 10).Value = 350
$ws.Cells.Item(367, 11).Value = 9000
$ws.Cells.Item(367, 12).Value = 10000
$ws.Cells.Item(367, 13).Value = 9343
$ws.Cells.Item(367, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(367, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(367, 16).Value = 187
$ws.Cells.Item(367, 17).Value = 50
$ws.Cells.Item(367, 18).Value = "Hortaliza"
